$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The two all-absent days (4th & 5th March, originally rows 7 & 8 with no
# attendance recorded) are removed from the tracker. This shifts the
# following rows up by two, so 6th March (old row 9) becomes row 7 and
# 7th March (old row 10) becomes row 8, and the two now-unused blank rows
# at the bottom (66 & 67) disappear.
$ws.Rows("7:8").Delete()

# Mark attendance (present) for some members on 27th Feb (row 4) and
# 28th Feb (row 5) for Paras, Hazel, George and Sam.
$ws.Range("C4:F4").Value = $true
$ws.Range("C5:F5").Value = $true

# Record attendance for 6th March (now row 7): Ben, Hazel, George and Sam
# present; Paras absent.
$ws.Range("B7").Value = $true
$ws.Range("D7:F7").Value = $true

# Leave the selection where the last edit was made.
$ws.Range("F4").Select() | Out-Null
